$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Leading apostrophe forces Excel to treat the value as literal text
    # (preventing auto-conversion of numeric-looking strings, which would
    # otherwise drop meaningful trailing zeros / precision), then reset the
    # cell style back to 'Normal' so no stray quote-prefix formatting sticks.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '90.997.53'
Set-TextValue $ws.Range('E2') '  +4.22%  '
Set-TextValue $ws.Range('D3') '3.097.86'
Set-TextValue $ws.Range('E3') '  +1.39%  '
Set-TextValue $ws.Range('E4') '  +0.26%  '
Set-TextValue $ws.Range('D5') '219.29'
Set-TextValue $ws.Range('E5') '  +4.94%  '
Set-TextValue $ws.Range('D6') '619.18'
Set-TextValue $ws.Range('E6') '  -0.29%  '
Set-TextValue $ws.Range('D7') '0.378'
Set-TextValue $ws.Range('E7') '  +5.15%  '
Set-TextValue $ws.Range('D8') '0.891'
Set-TextValue $ws.Range('E8') '  +15.92%  '
Set-TextValue $ws.Range('E9') '  +0.09%  '
Set-TextValue $ws.Range('D10') '3.093.76'
Set-TextValue $ws.Range('E10') '  +1.41%  '
Set-TextValue $ws.Range('D11') '0.695'
Set-TextValue $ws.Range('E11') '  +20.59%  '
Set-TextValue $ws.Range('E12') '  +7.48%  '
Set-TextValue $ws.Range('D13') '0.0000254'
Set-TextValue $ws.Range('E13') '  +9.21%  '
Set-TextValue $ws.Range('D14') '5.39'
Set-TextValue $ws.Range('E14') '  +3.34%  '
Set-TextValue $ws.Range('D15') '91.012.50'
Set-TextValue $ws.Range('E15') '  +4.52%  '
Set-TextValue $ws.Range('D16') '33.03'
Set-TextValue $ws.Range('E16') '  +6.12%  '
Set-TextValue $ws.Range('D17') '3.674.54'
Set-TextValue $ws.Range('E17') '  +1.85%  '
Set-TextValue $ws.Range('D18') '3.081.76'
Set-TextValue $ws.Range('E18') '  +0.81%  '
Set-TextValue $ws.Range('D19') '3.61'
Set-TextValue $ws.Range('E19') '  +7.90%  '
Set-TextValue $ws.Range('D20') '0.0000235'
Set-TextValue $ws.Range('E20') '  +13.42%  '
Set-TextValue $ws.Range('D21') '13.82'
Set-TextValue $ws.Range('E21') '  +6.81%  '
Set-TextValue $ws.Range('D22') '431.09'
Set-TextValue $ws.Range('E22') '  +4.20%  '
Set-TextValue $ws.Range('D23') '8.58'
Set-TextValue $ws.Range('E23') '  +5.38%  '
Set-TextValue $ws.Range('E24') '  +8.15%  '
Set-TextValue $ws.Range('D25') '5.58'
Set-TextValue $ws.Range('E25') '  +3.87%  '
Set-TextValue $ws.Range('D26') '11.90'
Set-TextValue $ws.Range('E26') '  +6.64%  '
Set-TextValue $ws.Range('D27') '83.68'
Set-TextValue $ws.Range('E27') '  +2.31%  '
Set-TextValue $ws.Range('D28') '3.263.25'
Set-TextValue $ws.Range('E28') '  +1.47%  '
Set-TextValue $ws.Range('E29') '  -0.19%  '
Set-TextValue $ws.Range('E30') '  +12.80%  '
Set-TextValue $ws.Range('E31') '  -0.04%  '
Set-TextValue $ws.Range('E32') '  +9.32%  '
Set-TextValue $ws.Range('D33') '3.93'
Set-TextValue $ws.Range('E33') '  +9.17%  '
Set-TextValue $ws.Range('D34') '517.83'
Set-TextValue $ws.Range('E34') '  +5.15%  '
Set-TextValue $ws.Range('D35') '6.94'
Set-TextValue $ws.Range('E35') '  +5.58%  '
Set-TextValue $ws.Range('D36') '0.140'
Set-TextValue $ws.Range('E36') '  +0.81%  '
Set-TextValue $ws.Range('D37') '1.29'
Set-TextValue $ws.Range('E37') '  +4.81%  '
Set-TextValue $ws.Range('E38') '  +3.57%  '
Set-TextValue $ws.Range('D39') '23.00'
Set-TextValue $ws.Range('E39') '  +5.86%  '
Set-TextValue $ws.Range('E40') '  +0.77%  '
Set-TextValue $ws.Range('E41') '  +0.11%  '
Set-TextValue $ws.Range('D42') '0.146'
Set-TextValue $ws.Range('E42') '  +12.40%  '
Set-TextValue $ws.Range('E43') '  +0.04%  '
Set-TextValue $ws.Range('D44') '0.372'
Set-TextValue $ws.Range('E44') '  +3.96%  '
Set-TextValue $ws.Range('E45') '  +4.28%  '
Set-TextValue $ws.Range('D46') '0.0718'
Set-TextValue $ws.Range('E46') '  +13.00%  '
Set-TextValue $ws.Range('D47') '43.87'
Set-TextValue $ws.Range('E47') '  +1.30%  '
Set-TextValue $ws.Range('D48') '141.47'
Set-TextValue $ws.Range('E48') '  -3.15%  '
Set-TextValue $ws.Range('D49') '1.26'
Set-TextValue $ws.Range('E49') '  +9.01%  '
Set-TextValue $ws.Range('D50') '0.000260'
Set-TextValue $ws.Range('E50') '  +20.31%  '
Set-TextValue $ws.Range('D51') '4.20'
Set-TextValue $ws.Range('E51') '  +9.15%  '
